$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (code) and B (date) hold values that look like numbers/dates
# ("512660", "2025-06-13"). Force each target cell to Text format *before*
# assigning so Excel stores them as plain text, matching the source data,
# instead of auto-converting to a number / date serial.

# Row 2: 512660 军工ETF
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "512660"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2025-06-13"
$ws.Cells.Item(2, 3).Value = "军工ETF"
$ws.Cells.Item(2, 4).Value = 1.06
$ws.Cells.Item(2, 5).Value = 1.44
$ws.Cells.Item(2, 6).Value = 1.34
$ws.Cells.Item(2, 7).Value = 7.88
$ws.Cells.Item(2, 8).Value = 1.053
$ws.Cells.Item(2, 9).Value = 1.0482
$ws.Cells.Item(2, 10).Value = 1.045
$ws.Cells.Item(2, 11).Value = -0.0001499999999998725
$ws.Cells.Item(2, 12).Value = $false
$ws.Cells.Item(2, 13).Value = 0.002500000000000169
$ws.Cells.Item(2, 14).Value = $false

# Row 3: 518880 黄金ETF
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "518880"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2025-06-13"
$ws.Cells.Item(3, 3).Value = "黄金ETF"
$ws.Cells.Item(3, 4).Value = 7.59
$ws.Cells.Item(3, 5).Value = 1.21
$ws.Cells.Item(3, 6).Value = 28
$ws.Cells.Item(3, 7).Value = 67.65
$ws.Cells.Item(3, 8).Value = 7.472
$ws.Cells.Item(3, 9).Value = 7.4355
$ws.Cells.Item(3, 10).Value = 7.372
$ws.Cells.Item(3, 11).Value = 0.01290000000000013
$ws.Cells.Item(3, 12).Value = $false
$ws.Cells.Item(3, 13).Value = 0.0224000000000002
$ws.Cells.Item(3, 14).Value = $false

# Row 4: 510410 资源ETF
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "510410"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2025-06-13"
$ws.Cells.Item(4, 3).Value = "资源ETF"
$ws.Cells.Item(4, 4).Value = 1.25
$ws.Cells.Item(4, 5).Value = 0.56
$ws.Cells.Item(4, 6).Value = 3.46
$ws.Cells.Item(4, 7).Value = 0.04
$ws.Cells.Item(4, 8).Value = 1.228
$ws.Cells.Item(4, 9).Value = 1.2184
$ws.Cells.Item(4, 10).Value = 1.215
$ws.Cells.Item(4, 11).Value = 0.002299999999999969
$ws.Cells.Item(4, 12).Value = $true
$ws.Cells.Item(4, 13).Value = 0.005100000000000104
$ws.Cells.Item(4, 14).Value = $false

# Row 5: 561560 电力ETF
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "561560"
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "2025-06-13"
$ws.Cells.Item(5, 3).Value = "电力ETF"
$ws.Cells.Item(5, 4).Value = 1.17
$ws.Cells.Item(5, 5).Value = 0.51
$ws.Cells.Item(5, 6).Value = -2
$ws.Cells.Item(5, 7).Value = 0.62
$ws.Cells.Item(5, 8).Value = 1.17
$ws.Cells.Item(5, 9).Value = 1.1728
$ws.Cells.Item(5, 10).Value = 1.178
$ws.Cells.Item(5, 11).Value = -0.0002500000000003055
$ws.Cells.Item(5, 12).Value = $false
$ws.Cells.Item(5, 13).Value = -0.001199999999999868
$ws.Cells.Item(5, 14).Value = $false

# Row 6: 159666 交通运输ETF
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "159666"
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "2025-06-13"
$ws.Cells.Item(6, 3).Value = "交通运输ETF"
$ws.Cells.Item(6, 4).Value = 0.99
$ws.Cells.Item(6, 5).Value = 0.2
$ws.Cells.Item(6, 6).Value = -0.5
$ws.Cells.Item(6, 7).Value = 0.04
$ws.Cells.Item(6, 8).Value = 0.984
$ws.Cells.Item(6, 9).Value = 0.9863
$ws.Cells.Item(6, 10).Value = 0.985
$ws.Cells.Item(6, 11).Value = 0.0007499999999999174
$ws.Cells.Item(6, 12).Value = $true
$ws.Cells.Item(6, 13).Value = 0.0007999999999999119
$ws.Cells.Item(6, 14).Value = $false

# Row 7: 516020 化工ETF
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "516020"
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 2).Value = "2025-06-13"
$ws.Cells.Item(7, 3).Value = "化工ETF"
$ws.Cells.Item(7, 4).Value = 0.6
$ws.Cells.Item(7, 5).Value = 0.17
$ws.Cells.Item(7, 6).Value = -0.99
$ws.Cells.Item(7, 7).Value = 0.1
$ws.Cells.Item(7, 8).Value = 0.596
$ws.Cells.Item(7, 9).Value = 0.5915999999999999
$ws.Cells.Item(7, 10).Value = 0.593
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = $false
$ws.Cells.Item(7, 13).Value = 0.001199999999999868
$ws.Cells.Item(7, 14).Value = $false

# Row 8: 515900 央企创新驱动ETF
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "515900"
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "2025-06-13"
$ws.Cells.Item(8, 3).Value = "央企创新驱动ETF"
$ws.Cells.Item(8, 4).Value = 1.43
$ws.Cells.Item(8, 5).Value = 0.14
$ws.Cells.Item(8, 6).Value = -5.05
$ws.Cells.Item(8, 7).Value = 0.19
$ws.Cells.Item(8, 8).Value = 1.431
$ws.Cells.Item(8, 9).Value = 1.425
$ws.Cells.Item(8, 10).Value = 1.427
$ws.Cells.Item(8, 11).Value = -0.000299999999999967
$ws.Cells.Item(8, 12).Value = $false
$ws.Cells.Item(8, 13).Value = 0.001800000000000024
$ws.Cells.Item(8, 14).Value = $false

# Row 9: 159691 港股红利ETF
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "159691"
$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "2025-06-13"
$ws.Cells.Item(9, 3).Value = "港股红利ETF"
$ws.Cells.Item(9, 4).Value = 1.2
$ws.Cells.Item(9, 5).Value = 0.08
$ws.Cells.Item(9, 6).Value = 8.27
$ws.Cells.Item(9, 7).Value = 1.84
$ws.Cells.Item(9, 8).Value = 1.187
$ws.Cells.Item(9, 9).Value = 1.1757
$ws.Cells.Item(9, 10).Value = 1.162
$ws.Cells.Item(9, 11).Value = $null
$ws.Cells.Item(9, 12).Value = $true
$ws.Cells.Item(9, 13).Value = 0.004899999999999904
$ws.Cells.Item(9, 14).Value = $true

# Row 10: 511090 30年国债ETF
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "511090"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "2025-06-13"
$ws.Cells.Item(10, 3).Value = "30年国债ETF"
$ws.Cells.Item(10, 4).Value = 124.24
$ws.Cells.Item(10, 5).Value = 0.01
$ws.Cells.Item(10, 6).Value = 1.29
$ws.Cells.Item(10, 7).Value = 58.01
$ws.Cells.Item(10, 8).Value = 123.669
$ws.Cells.Item(10, 9).Value = 123.2838
$ws.Cells.Item(10, 10).Value = 123.246
$ws.Cells.Item(10, 11).Value = 0.04565000000000907
$ws.Cells.Item(10, 12).Value = $false
$ws.Cells.Item(10, 13).Value = 0.128399999999985
$ws.Cells.Item(10, 14).Value = $false

# Row 11: 159652 有色50ETF
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "159652"
$ws.Cells.Item(11, 2).NumberFormat = "@"
$ws.Cells.Item(11, 2).Value = "2025-06-13"
$ws.Cells.Item(11, 3).Value = "有色50ETF"
$ws.Cells.Item(11, 4).Value = 0.97
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 13.25
$ws.Cells.Item(11, 7).Value = 0.21
$ws.Cells.Item(11, 8).Value = 0.941
$ws.Cells.Item(11, 9).Value = 0.9280000000000002
$ws.Cells.Item(11, 10).Value = 0.925
$ws.Cells.Item(11, 11).Value = 0.002449999999999952
$ws.Cells.Item(11, 12).Value = $true
$ws.Cells.Item(11, 13).Value = 0.006000000000000227
$ws.Cells.Item(11, 14).Value = $false

# Row 12: 515220 煤炭ETF
$ws.Cells.Item(12, 1).NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "515220"
$ws.Cells.Item(12, 2).NumberFormat = "@"
$ws.Cells.Item(12, 2).Value = "2025-06-13"
$ws.Cells.Item(12, 3).Value = "煤炭ETF"
$ws.Cells.Item(12, 4).Value = 0.98
$ws.Cells.Item(12, 5).Value = -0.2
$ws.Cells.Item(12, 6).Value = -13.06
$ws.Cells.Item(12, 7).Value = 1.13
$ws.Cells.Item(12, 8).Value = 0.985
$ws.Cells.Item(12, 9).Value = 0.9859
$ws.Cells.Item(12, 10).Value = 0.989
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = $false
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = $false

# Row 13: 512760 芯片ETF
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "512760"
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "2025-06-13"
$ws.Cells.Item(13, 3).Value = "芯片ETF"
$ws.Cells.Item(13, 4).Value = 1.1
$ws.Cells.Item(13, 5).Value = -0.27
$ws.Cells.Item(13, 6).Value = -2.31
$ws.Cells.Item(13, 7).Value = 1.98
$ws.Cells.Item(13, 8).Value = 1.125
$ws.Cells.Item(13, 9).Value = 1.1176
$ws.Cells.Item(13, 10).Value = 1.125
$ws.Cells.Item(13, 11).Value = -0.002799999999999914
$ws.Cells.Item(13, 12).Value = $false
$ws.Cells.Item(13, 13).Value = -0.0008000000000001339
$ws.Cells.Item(13, 14).Value = $false

# Row 14: 560070 央企红利ETF基金
$ws.Cells.Item(14, 1).NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "560070"
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "2025-06-13"
$ws.Cells.Item(14, 3).Value = "央企红利ETF基金"
$ws.Cells.Item(14, 4).Value = 1.02
$ws.Cells.Item(14, 5).Value = -0.29
$ws.Cells.Item(14, 6).Value = -3.77
$ws.Cells.Item(14, 7).Value = 0.02
$ws.Cells.Item(14, 8).Value = 1.022
$ws.Cells.Item(14, 9).Value = 1.0229
$ws.Cells.Item(14, 10).Value = 1.024
$ws.Cells.Item(14, 11).Value = 0.0002500000000000835
$ws.Cells.Item(14, 12).Value = $true
$ws.Cells.Item(14, 13).Value = 0.000299999999999967
$ws.Cells.Item(14, 14).Value = $false

# Row 15: 512480 半导体ETF
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "512480"
$ws.Cells.Item(15, 2).NumberFormat = "@"
$ws.Cells.Item(15, 2).Value = "2025-06-13"
$ws.Cells.Item(15, 3).Value = "半导体ETF"
$ws.Cells.Item(15, 4).Value = 0.99
$ws.Cells.Item(15, 5).Value = -0.4
$ws.Cells.Item(15, 6).Value = -1.59
$ws.Cells.Item(15, 7).Value = 8.17
$ws.Cells.Item(15, 8).Value = 1.016
$ws.Cells.Item(15, 9).Value = 1.0089
$ws.Cells.Item(15, 10).Value = 1.015
$ws.Cells.Item(15, 11).Value = -0.002650000000000041
$ws.Cells.Item(15, 12).Value = $false
$ws.Cells.Item(15, 13).Value = -0.000700000000000145
$ws.Cells.Item(15, 14).Value = $false

# Row 16: 588000 科创50ETF
$ws.Cells.Item(16, 1).NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "588000"
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "2025-06-13"
$ws.Cells.Item(16, 3).Value = "科创50ETF"
$ws.Cells.Item(16, 4).Value = 1.02
$ws.Cells.Item(16, 5).Value = -0.49
$ws.Cells.Item(16, 6).Value = -2.2
$ws.Cells.Item(16, 7).Value = 23.13
$ws.Cells.Item(16, 8).Value = 1.041
$ws.Cells.Item(16, 9).Value = 1.037
$ws.Cells.Item(16, 10).Value = 1.04
$ws.Cells.Item(16, 11).Value = -0.002000000000000002
$ws.Cells.Item(16, 12).Value = $false
$ws.Cells.Item(16, 13).Value = -0.0008000000000001339
$ws.Cells.Item(16, 14).Value = $false

# Row 17: 510050 上证50ETF
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "510050"
$ws.Cells.Item(17, 2).NumberFormat = "@"
$ws.Cells.Item(17, 2).Value = "2025-06-13"
$ws.Cells.Item(17, 3).Value = "上证50ETF"
$ws.Cells.Item(17, 4).Value = 2.75
$ws.Cells.Item(17, 5).Value = -0.51
$ws.Cells.Item(17, 6).Value = 0.29
$ws.Cells.Item(17, 7).Value = 17.88
$ws.Cells.Item(17, 8).Value = 2.752
$ws.Cells.Item(17, 9).Value = 2.7511
$ws.Cells.Item(17, 10).Value = 2.767
$ws.Cells.Item(17, 11).Value = -0.0009500000000000064
$ws.Cells.Item(17, 12).Value = $false
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = $false

# Row 18: 515080 中证红利ETF
$ws.Cells.Item(18, 1).NumberFormat = "@"
$ws.Cells.Item(18, 1).Value = "515080"
$ws.Cells.Item(18, 2).NumberFormat = "@"
$ws.Cells.Item(18, 2).Value = "2025-06-13"
$ws.Cells.Item(18, 3).Value = "中证红利ETF"
$ws.Cells.Item(18, 4).Value = 1.53
$ws.Cells.Item(18, 5).Value = -0.52
$ws.Cells.Item(18, 6).Value = -1.21
$ws.Cells.Item(18, 7).Value = 1.63
$ws.Cells.Item(18, 8).Value = 1.528
$ws.Cells.Item(18, 9).Value = 1.5263
$ws.Cells.Item(18, 10).Value = 1.527
$ws.Cells.Item(18, 11).Value = 0.0002500000000003055
$ws.Cells.Item(18, 12).Value = $true
$ws.Cells.Item(18, 13).Value = 0.001600000000000046
$ws.Cells.Item(18, 14).Value = $false

# Row 19: 159781 科创创业ETF
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "159781"
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = "2025-06-13"
$ws.Cells.Item(19, 3).Value = "科创创业ETF"
$ws.Cells.Item(19, 4).Value = 0.53
$ws.Cells.Item(19, 5).Value = -0.56
$ws.Cells.Item(19, 6).Value = -6.02
$ws.Cells.Item(19, 7).Value = 0.42
$ws.Cells.Item(19, 8).Value = 0.535
$ws.Cells.Item(19, 9).Value = 0.5307999999999999
$ws.Cells.Item(19, 10).Value = 0.532
$ws.Cells.Item(19, 11).Value = -0.0005499999999999394
$ws.Cells.Item(19, 12).Value = $false
$ws.Cells.Item(19, 13).Value = 0.0009999999999998899
$ws.Cells.Item(19, 14).Value = $false

# Row 20: 515260 电子ETF
$ws.Cells.Item(20, 1).NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "515260"
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = "2025-06-13"
$ws.Cells.Item(20, 3).Value = "电子ETF"
$ws.Cells.Item(20, 4).Value = 0.84
$ws.Cells.Item(20, 5).Value = -0.59
$ws.Cells.Item(20, 6).Value = -6.95
$ws.Cells.Item(20, 7).Value = 0.04
$ws.Cells.Item(20, 8).Value = 0.854
$ws.Cells.Item(20, 9).Value = 0.8440999999999999
$ws.Cells.Item(20, 10).Value = 0.847
$ws.Cells.Item(20, 11).Value = -0.00130000000000019
$ws.Cells.Item(20, 12).Value = $false
$ws.Cells.Item(20, 13).Value = 0.0007999999999999119
$ws.Cells.Item(20, 14).Value = $false

# Row 21: 512890 红利低波ETF
$ws.Cells.Item(21, 1).NumberFormat = "@"
$ws.Cells.Item(21, 1).Value = "512890"
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = "2025-06-13"
$ws.Cells.Item(21, 3).Value = "红利低波ETF"
$ws.Cells.Item(21, 4).Value = 1.17
$ws.Cells.Item(21, 5).Value = -0.6
$ws.Cells.Item(21, 6).Value = 4
$ws.Cells.Item(21, 7).Value = 9.3
$ws.Cells.Item(21, 8).Value = 1.168
$ws.Cells.Item(21, 9).Value = 1.1649
$ws.Cells.Item(21, 10).Value = 1.157
$ws.Cells.Item(21, 11).Value = 0.001299999999999857
$ws.Cells.Item(21, 12).Value = $true
$ws.Cells.Item(21, 13).Value = 0.002099999999999769
$ws.Cells.Item(21, 14).Value = $true

# Row 22: 513290 纳指生物科技ETF
$ws.Cells.Item(22, 1).NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = "513290"
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = "2025-06-13"
$ws.Cells.Item(22, 3).Value = "纳指生物科技ETF"
$ws.Cells.Item(22, 4).Value = 1.11
$ws.Cells.Item(22, 5).Value = -0.63
$ws.Cells.Item(22, 6).Value = -3.72
$ws.Cells.Item(22, 7).Value = 0.83
$ws.Cells.Item(22, 8).Value = 1.112
$ws.Cells.Item(22, 9).Value = 1.1049
$ws.Cells.Item(22, 10).Value = 1.093
$ws.Cells.Item(22, 11).Value = 0.002050000000000107
$ws.Cells.Item(22, 12).Value = $false
$ws.Cells.Item(22, 13).Value = 0.001499999999999835
$ws.Cells.Item(22, 14).Value = $false

# Row 23: 510300 沪深300ETF
$ws.Cells.Item(23, 1).NumberFormat = "@"
$ws.Cells.Item(23, 1).Value = "510300"
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = "2025-06-13"
$ws.Cells.Item(23, 3).Value = "沪深300ETF"
$ws.Cells.Item(23, 4).Value = 3.98
$ws.Cells.Item(23, 5).Value = -0.65
$ws.Cells.Item(23, 6).Value = -1.04
$ws.Cells.Item(23, 7).Value = 27.8
$ws.Cells.Item(23, 8).Value = 3.983
$ws.Cells.Item(23, 9).Value = 3.9715
$ws.Cells.Item(23, 10).Value = 3.984
$ws.Cells.Item(23, 11).Value = -0.0006999999999997009
$ws.Cells.Item(23, 12).Value = $false
$ws.Cells.Item(23, 13).Value = 0.003399999999999626
$ws.Cells.Item(23, 14).Value = $false

# Row 24: 159825 农业ETF
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "159825"
$ws.Cells.Item(24, 2).NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = "2025-06-13"
$ws.Cells.Item(24, 3).Value = "农业ETF"
$ws.Cells.Item(24, 4).Value = 0.72
$ws.Cells.Item(24, 5).Value = -0.69
$ws.Cells.Item(24, 6).Value = 8.56
$ws.Cells.Item(24, 7).Value = 0.48
$ws.Cells.Item(24, 8).Value = 0.716
$ws.Cells.Item(24, 9).Value = 0.708
$ws.Cells.Item(24, 10).Value = 0.698
$ws.Cells.Item(24, 11).Value = 0.001900000000000013
$ws.Cells.Item(24, 12).Value = $false
$ws.Cells.Item(24, 13).Value = 0.003099999999999992
$ws.Cells.Item(24, 14).Value = $false

# Row 25: 510210 上证指数ETF
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "510210"
$ws.Cells.Item(25, 2).NumberFormat = "@"
$ws.Cells.Item(25, 2).Value = "2025-06-13"
$ws.Cells.Item(25, 3).Value = "上证指数ETF"
$ws.Cells.Item(25, 4).Value = 0.84
$ws.Cells.Item(25, 5).Value = -0.71
$ws.Cells.Item(25, 6).Value = 3.06
$ws.Cells.Item(25, 7).Value = 1.23
$ws.Cells.Item(25, 8).Value = 0.841
$ws.Cells.Item(25, 9).Value = 0.8368
$ws.Cells.Item(25, 10).Value = 0.833
$ws.Cells.Item(25, 11).Value = 0.0007999999999999119
$ws.Cells.Item(25, 12).Value = $true
$ws.Cells.Item(25, 13).Value = 0.001599999999999935
$ws.Cells.Item(25, 14).Value = $false

# Row 26: 510760 上证综指ETF
$ws.Cells.Item(26, 1).NumberFormat = "@"
$ws.Cells.Item(26, 1).Value = "510760"
$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = "2025-06-13"
$ws.Cells.Item(26, 3).Value = "上证综指ETF"
$ws.Cells.Item(26, 4).Value = 1.08
$ws.Cells.Item(26, 5).Value = -0.73
$ws.Cells.Item(26, 6).Value = 2.85
$ws.Cells.Item(26, 7).Value = 0.43
$ws.Cells.Item(26, 8).Value = 1.081
$ws.Cells.Item(26, 9).Value = 1.0763
$ws.Cells.Item(26, 10).Value = 1.074
$ws.Cells.Item(26, 11).Value = 0.0005999999999999339
$ws.Cells.Item(26, 12).Value = $true
$ws.Cells.Item(26, 13).Value = 0.001600000000000046
$ws.Cells.Item(26, 14).Value = $false

# Row 27: 513520 日经ETF
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "513520"
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = "2025-06-13"
$ws.Cells.Item(27, 3).Value = "日经ETF"
$ws.Cells.Item(27, 4).Value = 1.46
$ws.Cells.Item(27, 5).Value = -0.75
$ws.Cells.Item(27, 6).Value = 0.41
$ws.Cells.Item(27, 7).Value = 0.53
$ws.Cells.Item(27, 8).Value = 1.461
$ws.Cells.Item(27, 9).Value = 1.4613
$ws.Cells.Item(27, 10).Value = 1.455
$ws.Cells.Item(27, 11).Value = 0.001100000000000101
$ws.Cells.Item(27, 12).Value = $true
$ws.Cells.Item(27, 13).Value = -0.0007999999999999119
$ws.Cells.Item(27, 14).Value = $false

# Row 28: 512800 银行ETF
$ws.Cells.Item(28, 1).NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "512800"
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = "2025-06-13"
$ws.Cells.Item(28, 3).Value = "银行ETF"
$ws.Cells.Item(28, 4).Value = 1.65
$ws.Cells.Item(28, 5).Value = -0.78
$ws.Cells.Item(28, 6).Value = 11.11
$ws.Cells.Item(28, 7).Value = 7.94
$ws.Cells.Item(28, 8).Value = 1.643
$ws.Cells.Item(28, 9).Value = 1.6327
$ws.Cells.Item(28, 10).Value = 1.621
$ws.Cells.Item(28, 11).Value = 0.002249999999999863
$ws.Cells.Item(28, 12).Value = $true
$ws.Cells.Item(28, 13).Value = 0.004299999999999971
$ws.Cells.Item(28, 14).Value = $false

# Row 29: 515800 800ETF
$ws.Cells.Item(29, 1).NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "515800"
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = "2025-06-13"
$ws.Cells.Item(29, 3).Value = "800ETF"
$ws.Cells.Item(29, 4).Value = 1.01
$ws.Cells.Item(29, 5).Value = -0.79
$ws.Cells.Item(29, 6).Value = -0.49
$ws.Cells.Item(29, 7).Value = 0.72
$ws.Cells.Item(29, 8).Value = 1.009
$ws.Cells.Item(29, 9).Value = 1.0041
$ws.Cells.Item(29, 10).Value = 1.006
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = $false
$ws.Cells.Item(29, 13).Value = 0.001500000000000057
$ws.Cells.Item(29, 14).Value = $false

# Row 30: 510230 金融ETF
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "510230"
$ws.Cells.Item(30, 2).NumberFormat = "@"
$ws.Cells.Item(30, 2).Value = "2025-06-13"
$ws.Cells.Item(30, 3).Value = "金融ETF"
$ws.Cells.Item(30, 4).Value = 1.38
$ws.Cells.Item(30, 5).Value = -0.79
$ws.Cells.Item(30, 6).Value = 5.89
$ws.Cells.Item(30, 7).Value = 0.19
$ws.Cells.Item(30, 8).Value = 1.376
$ws.Cells.Item(30, 9).Value = 1.3656
$ws.Cells.Item(30, 10).Value = 1.36
$ws.Cells.Item(30, 11).Value = 0.001649999999999707
$ws.Cells.Item(30, 12).Value = $true
$ws.Cells.Item(30, 13).Value = 0.003500000000000281
$ws.Cells.Item(30, 14).Value = $false

# Row 31: 513800 日本东证指数ETF
$ws.Cells.Item(31, 1).NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "513800"
$ws.Cells.Item(31, 2).NumberFormat = "@"
$ws.Cells.Item(31, 2).Value = "2025-06-13"
$ws.Cells.Item(31, 3).Value = "日本东证指数ETF"
$ws.Cells.Item(31, 4).Value = 1.47
$ws.Cells.Item(31, 5).Value = -0.81
$ws.Cells.Item(31, 6).Value = 7.39
$ws.Cells.Item(31, 7).Value = 0.14
$ws.Cells.Item(31, 8).Value = 1.475
$ws.Cells.Item(31, 9).Value = 1.4802
$ws.Cells.Item(31, 10).Value = 1.47
$ws.Cells.Item(31, 11).Value = 0.001400000000000068
$ws.Cells.Item(31, 12).Value = $true
$ws.Cells.Item(31, 13).Value = -0.001600000000000046
$ws.Cells.Item(31, 14).Value = $false

# Row 32: 512500 中证500ETF华夏
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "512500"
$ws.Cells.Item(32, 2).NumberFormat = "@"
$ws.Cells.Item(32, 2).Value = "2025-06-13"
$ws.Cells.Item(32, 3).Value = "中证500ETF华夏"
$ws.Cells.Item(32, 4).Value = 3.17
$ws.Cells.Item(32, 5).Value = -0.81
$ws.Cells.Item(32, 6).Value = 0.96
$ws.Cells.Item(32, 7).Value = 1.94
$ws.Cells.Item(32, 8).Value = 3.176
$ws.Cells.Item(32, 9).Value = 3.1533
$ws.Cells.Item(32, 10).Value = 3.144
$ws.Cells.Item(32, 11).Value = 0.0000499999999998834709913
$ws.Cells.Item(32, 12).Value = $false
$ws.Cells.Item(32, 13).Value = 0.006600000000000161
$ws.Cells.Item(32, 14).Value = $false

# Row 33: 515210 钢铁ETF
$ws.Cells.Item(33, 1).NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "515210"
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = "2025-06-13"
$ws.Cells.Item(33, 3).Value = "钢铁ETF"
$ws.Cells.Item(33, 4).Value = 1.2
$ws.Cells.Item(33, 5).Value = -0.83
$ws.Cells.Item(33, 6).Value = 2.3
$ws.Cells.Item(33, 7).Value = 0.34
$ws.Cells.Item(33, 8).Value = 1.202
$ws.Cells.Item(33, 9).Value = 1.1996
$ws.Cells.Item(33, 10).Value = 1.211
$ws.Cells.Item(33, 11).Value = -0.001650000000000151
$ws.Cells.Item(33, 12).Value = $false
$ws.Cells.Item(33, 13).Value = -0.000299999999999967
$ws.Cells.Item(33, 14).Value = $false

# Row 34: 513500 标普500ETF
$ws.Cells.Item(34, 1).NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "513500"
$ws.Cells.Item(34, 2).NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = "2025-06-13"
$ws.Cells.Item(34, 3).Value = "标普500ETF"
$ws.Cells.Item(34, 4).Value = 2.02
$ws.Cells.Item(34, 5).Value = -0.83
$ws.Cells.Item(34, 6).Value = -6.43
$ws.Cells.Item(34, 7).Value = 4.49
$ws.Cells.Item(34, 8).Value = 2.028
$ws.Cells.Item(34, 9).Value = 2.0227
$ws.Cells.Item(34, 10).Value = 2.013
$ws.Cells.Item(34, 11).Value = 0.002150000000000318
$ws.Cells.Item(34, 12).Value = $true
$ws.Cells.Item(34, 13).Value = 0.002299999999999969
$ws.Cells.Item(34, 14).Value = $false

# Row 35: 159949 创业板50ETF
$ws.Cells.Item(35, 1).NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "159949"
$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = "2025-06-13"
$ws.Cells.Item(35, 3).Value = "创业板50ETF"
$ws.Cells.Item(35, 4).Value = 0.91
$ws.Cells.Item(35, 5).Value = -0.88
$ws.Cells.Item(35, 6).Value = -5.72
$ws.Cells.Item(35, 7).Value = 5.84
$ws.Cells.Item(35, 8).Value = 0.906
$ws.Cells.Item(35, 9).Value = 0.8956
$ws.Cells.Item(35, 10).Value = 0.9
$ws.Cells.Item(35, 11).Value = -0.000200000000000089
$ws.Cells.Item(35, 12).Value = $false
$ws.Cells.Item(35, 13).Value = 0.002599999999999936
$ws.Cells.Item(35, 14).Value = $false

# Row 36: 515880 通信ETF
$ws.Cells.Item(36, 1).NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "515880"
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = "2025-06-13"
$ws.Cells.Item(36, 3).Value = "通信ETF"
$ws.Cells.Item(36, 4).Value = 1.32
$ws.Cells.Item(36, 5).Value = -0.9
$ws.Cells.Item(36, 6).Value = -3.08
$ws.Cells.Item(36, 7).Value = 1.79
$ws.Cells.Item(36, 8).Value = 1.312
$ws.Cells.Item(36, 9).Value = 1.2748
$ws.Cells.Item(36, 10).Value = 1.253
$ws.Cells.Item(36, 11).Value = 0.003149999999999986
$ws.Cells.Item(36, 12).Value = $true
$ws.Cells.Item(36, 13).Value = 0.01190000000000024
$ws.Cells.Item(36, 14).Value = $false

# Row 37: 159915 创业板ETF
$ws.Cells.Item(37, 1).NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "159915"
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = "2025-06-13"
$ws.Cells.Item(37, 3).Value = "创业板ETF"
$ws.Cells.Item(37, 4).Value = 2.02
$ws.Cells.Item(37, 5).Value = -0.98
$ws.Cells.Item(37, 6).Value = -3.9
$ws.Cells.Item(37, 7).Value = 14.8
$ws.Cells.Item(37, 8).Value = 2.021
$ws.Cells.Item(37, 9).Value = 2.0003
$ws.Cells.Item(37, 10).Value = 2.004
$ws.Cells.Item(37, 11).Value = -0.000500000000000167
$ws.Cells.Item(37, 12).Value = $false
$ws.Cells.Item(37, 13).Value = 0.005300000000000304
$ws.Cells.Item(37, 14).Value = $false

# Row 38: 516670 畜牧养殖ETF
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "516670"
$ws.Cells.Item(38, 2).NumberFormat = "@"
$ws.Cells.Item(38, 2).Value = "2025-06-13"
$ws.Cells.Item(38, 3).Value = "畜牧养殖ETF"
$ws.Cells.Item(38, 4).Value = 0.67
$ws.Cells.Item(38, 5).Value = -1.04
$ws.Cells.Item(38, 6).Value = 7.09
$ws.Cells.Item(38, 7).Value = 0.13
$ws.Cells.Item(38, 8).Value = 0.665
$ws.Cells.Item(38, 9).Value = 0.6612
$ws.Cells.Item(38, 10).Value = 0.656
$ws.Cells.Item(38, 11).Value = 0.001449999999999951
$ws.Cells.Item(38, 12).Value = $false
$ws.Cells.Item(38, 13).Value = 0.001900000000000013
$ws.Cells.Item(38, 14).Value = $false

# Row 39: 512000 券商ETF
$ws.Cells.Item(39, 1).NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "512000"
$ws.Cells.Item(39, 2).NumberFormat = "@"
$ws.Cells.Item(39, 2).Value = "2025-06-13"
$ws.Cells.Item(39, 3).Value = "券商ETF"
$ws.Cells.Item(39, 4).Value = 1.03
$ws.Cells.Item(39, 5).Value = -1.06
$ws.Cells.Item(39, 6).Value = -8.13
$ws.Cells.Item(39, 7).Value = 6.18
$ws.Cells.Item(39, 8).Value = 1.027
$ws.Cells.Item(39, 9).Value = 1.015
$ws.Cells.Item(39, 10).Value = 1.016
$ws.Cells.Item(39, 11).Value = 0.0007500000000002505
$ws.Cells.Item(39, 12).Value = $false
$ws.Cells.Item(39, 13).Value = 0.003500000000000059
$ws.Cells.Item(39, 14).Value = $false

# Row 40: 513100 纳指ETF
$ws.Cells.Item(40, 1).NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = "513100"
$ws.Cells.Item(40, 2).NumberFormat = "@"
$ws.Cells.Item(40, 2).Value = "2025-06-13"
$ws.Cells.Item(40, 3).Value = "纳指ETF"
$ws.Cells.Item(40, 4).Value = 1.55
$ws.Cells.Item(40, 5).Value = -1.08
$ws.Cells.Item(40, 6).Value = -4.6
$ws.Cells.Item(40, 7).Value = 9.97
$ws.Cells.Item(40, 8).Value = 1.566
$ws.Cells.Item(40, 9).Value = 1.568
$ws.Cells.Item(40, 10).Value = 1.56
$ws.Cells.Item(40, 11).Value = 0.00140000000000029
$ws.Cells.Item(40, 12).Value = $true
$ws.Cells.Item(40, 13).Value = 0
$ws.Cells.Item(40, 14).Value = $false

# Row 41: 515070 人工智能AIETF
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = "515070"
$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "2025-06-13"
$ws.Cells.Item(41, 3).Value = "人工智能AIETF"
$ws.Cells.Item(41, 4).Value = 1.16
$ws.Cells.Item(41, 5).Value = -1.28
$ws.Cells.Item(41, 6).Value = -0.09
$ws.Cells.Item(41, 7).Value = 0.84
$ws.Cells.Item(41, 8).Value = 1.181
$ws.Cells.Item(41, 9).Value = 1.1638
$ws.Cells.Item(41, 10).Value = 1.164
$ws.Cells.Item(41, 11).Value = -0.001449999999999951
$ws.Cells.Item(41, 12).Value = $false
$ws.Cells.Item(41, 13).Value = 0.002499999999999947
$ws.Cells.Item(41, 14).Value = $false

# Row 42: 159637 新能源车龙头ETF
$ws.Cells.Item(42, 1).NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = "159637"
$ws.Cells.Item(42, 2).NumberFormat = "@"
$ws.Cells.Item(42, 2).Value = "2025-06-13"
$ws.Cells.Item(42, 3).Value = "新能源车龙头ETF"
$ws.Cells.Item(42, 4).Value = 0.6
$ws.Cells.Item(42, 5).Value = -1.32
$ws.Cells.Item(42, 6).Value = 1.19
$ws.Cells.Item(42, 7).Value = 0.09
$ws.Cells.Item(42, 8).Value = 0.6
$ws.Cells.Item(42, 9).Value = 0.5986
$ws.Cells.Item(42, 10).Value = 0.608
$ws.Cells.Item(42, 11).Value = -0.0009000000000000119
$ws.Cells.Item(42, 12).Value = $false
$ws.Cells.Item(42, 13).Value = -0.000299999999999967
$ws.Cells.Item(42, 14).Value = $false

# Row 43: 515230 软件ETF
$ws.Cells.Item(43, 1).NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = "515230"
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = "2025-06-13"
$ws.Cells.Item(43, 3).Value = "软件ETF"
$ws.Cells.Item(43, 4).Value = 0.78
$ws.Cells.Item(43, 5).Value = -1.39
$ws.Cells.Item(43, 6).Value = -0.13
$ws.Cells.Item(43, 7).Value = 0.52
$ws.Cells.Item(43, 8).Value = 0.801
$ws.Cells.Item(43, 9).Value = 0.7925000000000001
$ws.Cells.Item(43, 10).Value = 0.793
$ws.Cells.Item(43, 11).Value = -0.001900000000000124
$ws.Cells.Item(43, 12).Value = $false
$ws.Cells.Item(43, 13).Value = 0.0006000000000000449
$ws.Cells.Item(43, 14).Value = $false

# Row 44: 159667 工业母机ETF
$ws.Cells.Item(44, 1).NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = "159667"
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = "2025-06-13"
$ws.Cells.Item(44, 3).Value = "工业母机ETF"
$ws.Cells.Item(44, 4).Value = 1.11
$ws.Cells.Item(44, 5).Value = -1.42
$ws.Cells.Item(44, 6).Value = 9.5
$ws.Cells.Item(44, 7).Value = 0.55
$ws.Cells.Item(44, 8).Value = 1.118
$ws.Cells.Item(44, 9).Value = 1.1093
$ws.Cells.Item(44, 10).Value = 1.117
$ws.Cells.Item(44, 11).Value = -0.001800000000000246
$ws.Cells.Item(44, 12).Value = $false
$ws.Cells.Item(44, 13).Value = 0.001199999999999868
$ws.Cells.Item(44, 14).Value = $false

# Row 45: 512170 医疗ETF
$ws.Cells.Item(45, 1).NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = "512170"
$ws.Cells.Item(45, 2).NumberFormat = "@"
$ws.Cells.Item(45, 2).Value = "2025-06-13"
$ws.Cells.Item(45, 3).Value = "医疗ETF"
$ws.Cells.Item(45, 4).Value = 0.33
$ws.Cells.Item(45, 5).Value = -1.49
$ws.Cells.Item(45, 6).Value = 1.22
$ws.Cells.Item(45, 7).Value = 5.94
$ws.Cells.Item(45, 8).Value = 0.332
$ws.Cells.Item(45, 9).Value = 0.3307
$ws.Cells.Item(45, 10).Value = 0.327
$ws.Cells.Item(45, 11).Value = 0.0002000000000000335
$ws.Cells.Item(45, 12).Value = $true
$ws.Cells.Item(45, 13).Value = 0.0006000000000000449
$ws.Cells.Item(45, 14).Value = $false

# Row 46: 159770 机器人ETF
$ws.Cells.Item(46, 1).NumberFormat = "@"
$ws.Cells.Item(46, 1).Value = "159770"
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = "2025-06-13"
$ws.Cells.Item(46, 3).Value = "机器人ETF"
$ws.Cells.Item(46, 4).Value = 0.85
$ws.Cells.Item(46, 5).Value = -1.51
$ws.Cells.Item(46, 6).Value = 5.33
$ws.Cells.Item(46, 7).Value = 1.53
$ws.Cells.Item(46, 8).Value = 0.866
$ws.Cells.Item(46, 9).Value = 0.865
$ws.Cells.Item(46, 10).Value = 0.879
$ws.Cells.Item(46, 11).Value = -0.003399999999999848
$ws.Cells.Item(46, 12).Value = $false
$ws.Cells.Item(46, 13).Value = -0.001500000000000057
$ws.Cells.Item(46, 14).Value = $false

# Row 47: 515790 光伏ETF
$ws.Cells.Item(47, 1).NumberFormat = "@"
$ws.Cells.Item(47, 1).Value = "515790"
$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = "2025-06-13"
$ws.Cells.Item(47, 3).Value = "光伏ETF"
$ws.Cells.Item(47, 4).Value = 0.64
$ws.Cells.Item(47, 5).Value = -1.53
$ws.Cells.Item(47, 6).Value = -15.06
$ws.Cells.Item(47, 7).Value = 2.13
$ws.Cells.Item(47, 8).Value = 0.648
$ws.Cells.Item(47, 9).Value = 0.6452
$ws.Cells.Item(47, 10).Value = 0.654
$ws.Cells.Item(47, 11).Value = -0.002349999999999963
$ws.Cells.Item(47, 12).Value = $false
$ws.Cells.Item(47, 13).Value = -0.0007000000000000339
$ws.Cells.Item(47, 14).Value = $false

# Row 48: 512200 房地产ETF
$ws.Cells.Item(48, 1).NumberFormat = "@"
$ws.Cells.Item(48, 1).Value = "512200"
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = "2025-06-13"
$ws.Cells.Item(48, 3).Value = "房地产ETF"
$ws.Cells.Item(48, 4).Value = 1.33
$ws.Cells.Item(48, 5).Value = -1.55
$ws.Cells.Item(48, 6).Value = -8.89
$ws.Cells.Item(48, 7).Value = 1.34
$ws.Cells.Item(48, 8).Value = 1.35
$ws.Cells.Item(48, 9).Value = 1.3467
$ws.Cells.Item(48, 10).Value = 1.35
$ws.Cells.Item(48, 11).Value = -0.00175000000000014
$ws.Cells.Item(48, 12).Value = $false
$ws.Cells.Item(48, 13).Value = -0.0007999999999999119
$ws.Cells.Item(48, 14).Value = $false

# Row 49: 512010 医药ETF
$ws.Cells.Item(49, 1).NumberFormat = "@"
$ws.Cells.Item(49, 1).Value = "512010"
$ws.Cells.Item(49, 2).NumberFormat = "@"
$ws.Cells.Item(49, 2).Value = "2025-06-13"
$ws.Cells.Item(49, 3).Value = "医药ETF"
$ws.Cells.Item(49, 4).Value = 0.38
$ws.Cells.Item(49, 5).Value = -1.57
$ws.Cells.Item(49, 6).Value = 4.43
$ws.Cells.Item(49, 7).Value = 7.41
$ws.Cells.Item(49, 8).Value = 0.377
$ws.Cells.Item(49, 9).Value = 0.3752
$ws.Cells.Item(49, 10).Value = 0.371
$ws.Cells.Item(49, 11).Value = 0.0009999999999999454
$ws.Cells.Item(49, 12).Value = $true
$ws.Cells.Item(49, 13).Value = 0.001299999999999968
$ws.Cells.Item(49, 14).Value = $true

# Row 50: 510900 H股ETF
$ws.Cells.Item(50, 1).NumberFormat = "@"
$ws.Cells.Item(50, 1).Value = "510900"
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = "2025-06-13"
$ws.Cells.Item(50, 3).Value = "H股ETF"
$ws.Cells.Item(50, 4).Value = 1.11
$ws.Cells.Item(50, 5).Value = -1.59
$ws.Cells.Item(50, 6).Value = 17.26
$ws.Cells.Item(50, 7).Value = 2.95
$ws.Cells.Item(50, 8).Value = 1.119
$ws.Cells.Item(50, 9).Value = 1.107
$ws.Cells.Item(50, 10).Value = 1.105
$ws.Cells.Item(50, 11).Value = 0.001149999999999984
$ws.Cells.Item(50, 12).Value = $true
$ws.Cells.Item(50, 13).Value = 0.001900000000000013
$ws.Cells.Item(50, 14).Value = $false

# Row 51: 562390 中药50ETF
$ws.Cells.Item(51, 1).NumberFormat = "@"
$ws.Cells.Item(51, 1).Value = "562390"
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = "2025-06-13"
$ws.Cells.Item(51, 3).Value = "中药50ETF"
$ws.Cells.Item(51, 4).Value = 0.98
$ws.Cells.Item(51, 5).Value = -1.61
$ws.Cells.Item(51, 6).Value = -2.2
$ws.Cells.Item(51, 7).Value = 0.02
$ws.Cells.Item(51, 8).Value = 0.984
$ws.Cells.Item(51, 9).Value = 0.9810000000000001
$ws.Cells.Item(51, 10).Value = 0.972
$ws.Cells.Item(51, 11).Value = 0.001399999999999957
$ws.Cells.Item(51, 12).Value = $true
$ws.Cells.Item(51, 13).Value = 0.001500000000000057
$ws.Cells.Item(51, 14).Value = $true

# Row 52: 516510 云计算ETF
$ws.Cells.Item(52, 1).NumberFormat = "@"
$ws.Cells.Item(52, 1).Value = "516510"
$ws.Cells.Item(52, 2).NumberFormat = "@"
$ws.Cells.Item(52, 2).Value = "2025-06-13"
$ws.Cells.Item(52, 3).Value = "云计算ETF"
$ws.Cells.Item(52, 4).Value = 1.11
$ws.Cells.Item(52, 5).Value = -1.77
$ws.Cells.Item(52, 6).Value = 0.82
$ws.Cells.Item(52, 7).Value = 1.62
$ws.Cells.Item(52, 8).Value = 1.139
$ws.Cells.Item(52, 9).Value = 1.1203
$ws.Cells.Item(52, 10).Value = 1.117
$ws.Cells.Item(52, 11).Value = -0.002250000000000085
$ws.Cells.Item(52, 12).Value = $false
$ws.Cells.Item(52, 13).Value = 0.00199999999999978
$ws.Cells.Item(52, 14).Value = $false

# Row 53: 516010 游戏ETF
$ws.Cells.Item(53, 1).NumberFormat = "@"
$ws.Cells.Item(53, 1).Value = "516010"
$ws.Cells.Item(53, 2).NumberFormat = "@"
$ws.Cells.Item(53, 2).Value = "2025-06-13"
$ws.Cells.Item(53, 3).Value = "游戏ETF"
$ws.Cells.Item(53, 4).Value = 1.18
$ws.Cells.Item(53, 5).Value = -1.91
$ws.Cells.Item(53, 6).Value = 18.66
$ws.Cells.Item(53, 7).Value = 1.37
$ws.Cells.Item(53, 8).Value = 1.163
$ws.Cells.Item(53, 9).Value = 1.1413
$ws.Cells.Item(53, 10).Value = 1.115
$ws.Cells.Item(53, 11).Value = 0.004149999999999876
$ws.Cells.Item(53, 12).Value = $true
$ws.Cells.Item(53, 13).Value = 0.00749999999999984
$ws.Cells.Item(53, 14).Value = $false

# Row 54: 159928 消费ETF
$ws.Cells.Item(54, 1).NumberFormat = "@"
$ws.Cells.Item(54, 1).Value = "159928"
$ws.Cells.Item(54, 2).NumberFormat = "@"
$ws.Cells.Item(54, 2).Value = "2025-06-13"
$ws.Cells.Item(54, 3).Value = "消费ETF"
$ws.Cells.Item(54, 4).Value = 0.8
$ws.Cells.Item(54, 5).Value = -1.97
$ws.Cells.Item(54, 6).Value = -2.33
$ws.Cells.Item(54, 7).Value = 3.35
$ws.Cells.Item(54, 8).Value = 0.816
$ws.Cells.Item(54, 9).Value = 0.8184999999999999
$ws.Cells.Item(54, 10).Value = 0.821
$ws.Cells.Item(54, 11).Value = -0.001149999999999984
$ws.Cells.Item(54, 12).Value = $false
$ws.Cells.Item(54, 13).Value = -0.001600000000000157
$ws.Cells.Item(54, 14).Value = $false

# Row 55: 515250 智能汽车ETF
$ws.Cells.Item(55, 1).NumberFormat = "@"
$ws.Cells.Item(55, 1).Value = "515250"
$ws.Cells.Item(55, 2).NumberFormat = "@"
$ws.Cells.Item(55, 2).Value = "2025-06-13"
$ws.Cells.Item(55, 3).Value = "智能汽车ETF"
$ws.Cells.Item(55, 4).Value = 0.93
$ws.Cells.Item(55, 5).Value = -2.01
$ws.Cells.Item(55, 6).Value = -0.75
$ws.Cells.Item(55, 7).Value = 0.25
$ws.Cells.Item(55, 8).Value = 0.948
$ws.Cells.Item(55, 9).Value = 0.9477
$ws.Cells.Item(55, 10).Value = 0.953
$ws.Cells.Item(55, 11).Value = -0.002050000000000107
$ws.Cells.Item(55, 12).Value = $false
$ws.Cells.Item(55, 13).Value = -0.0007000000000000339
$ws.Cells.Item(55, 14).Value = $false

# Row 56: 159883 医疗器械ETF
$ws.Cells.Item(56, 1).NumberFormat = "@"
$ws.Cells.Item(56, 1).Value = "159883"
$ws.Cells.Item(56, 2).NumberFormat = "@"
$ws.Cells.Item(56, 2).Value = "2025-06-13"
$ws.Cells.Item(56, 3).Value = "医疗器械ETF"
$ws.Cells.Item(56, 4).Value = 0.48
$ws.Cells.Item(56, 5).Value = -2.03
$ws.Cells.Item(56, 6).Value = -0.82
$ws.Cells.Item(56, 7).Value = 0.73
$ws.Cells.Item(56, 8).Value = 0.487
$ws.Cells.Item(56, 9).Value = 0.4848
$ws.Cells.Item(56, 10).Value = 0.48
$ws.Cells.Item(56, 11).Value = 0.0000999999999999334754364
$ws.Cells.Item(56, 12).Value = $true
$ws.Cells.Item(56, 13).Value = 0.0005000000000000004
$ws.Cells.Item(56, 14).Value = $false

# Row 57: 512980 传媒ETF
$ws.Cells.Item(57, 1).NumberFormat = "@"
$ws.Cells.Item(57, 1).Value = "512980"
$ws.Cells.Item(57, 2).NumberFormat = "@"
$ws.Cells.Item(57, 2).Value = "2025-06-13"
$ws.Cells.Item(57, 3).Value = "传媒ETF"
$ws.Cells.Item(57, 4).Value = 0.82
$ws.Cells.Item(57, 5).Value = -2.26
$ws.Cells.Item(57, 6).Value = 6.33
$ws.Cells.Item(57, 7).Value = 1.48
$ws.Cells.Item(57, 8).Value = 0.823
$ws.Cells.Item(57, 9).Value = 0.8141
$ws.Cells.Item(57, 10).Value = 0.806
$ws.Cells.Item(57, 11).Value = 0.0006999999999999229
$ws.Cells.Item(57, 12).Value = $false
$ws.Cells.Item(57, 13).Value = 0.00240000000000018
$ws.Cells.Item(57, 14).Value = $false

# Row 58: 159643 疫苗ETF
$ws.Cells.Item(58, 1).NumberFormat = "@"
$ws.Cells.Item(58, 1).Value = "159643"
$ws.Cells.Item(58, 2).NumberFormat = "@"
$ws.Cells.Item(58, 2).Value = "2025-06-13"
$ws.Cells.Item(58, 3).Value = "疫苗ETF"
$ws.Cells.Item(58, 4).Value = 0.59
$ws.Cells.Item(58, 5).Value = -2.33
$ws.Cells.Item(58, 6).Value = -0.68
$ws.Cells.Item(58, 7).Value = 0.07
$ws.Cells.Item(58, 8).Value = 0.592
$ws.Cells.Item(58, 9).Value = 0.585
$ws.Cells.Item(58, 10).Value = 0.576
$ws.Cells.Item(58, 11).Value = 0.001449999999999951
$ws.Cells.Item(58, 12).Value = $true
$ws.Cells.Item(58, 13).Value = 0.002000000000000002
$ws.Cells.Item(58, 14).Value = $true

# Row 59: 513060 恒生医疗ETF
$ws.Cells.Item(59, 1).NumberFormat = "@"
$ws.Cells.Item(59, 1).Value = "513060"
$ws.Cells.Item(59, 2).NumberFormat = "@"
$ws.Cells.Item(59, 2).Value = "2025-06-13"
$ws.Cells.Item(59, 3).Value = "恒生医疗ETF"
$ws.Cells.Item(59, 4).Value = 0.58
$ws.Cells.Item(59, 5).Value = -2.34
$ws.Cells.Item(59, 6).Value = 49.74
$ws.Cells.Item(59, 7).Value = 39.4
$ws.Cells.Item(59, 8).Value = 0.558
$ws.Cells.Item(59, 9).Value = 0.539
$ws.Cells.Item(59, 10).Value = 0.515
$ws.Cells.Item(59, 11).Value = 0.005400000000000071
$ws.Cells.Item(59, 12).Value = $true
$ws.Cells.Item(59, 13).Value = 0.007800000000000029
$ws.Cells.Item(59, 14).Value = $false

# Row 60: 513330 恒生互联网ETF
$ws.Cells.Item(60, 1).NumberFormat = "@"
$ws.Cells.Item(60, 1).Value = "513330"
$ws.Cells.Item(60, 2).NumberFormat = "@"
$ws.Cells.Item(60, 2).Value = "2025-06-13"
$ws.Cells.Item(60, 3).Value = "恒生互联网ETF"
$ws.Cells.Item(60, 4).Value = 0.48
$ws.Cells.Item(60, 5).Value = -2.42
$ws.Cells.Item(60, 6).Value = 16.07
$ws.Cells.Item(60, 7).Value = 22.54
$ws.Cells.Item(60, 8).Value = 0.487
$ws.Cells.Item(60, 9).Value = 0.4791
$ws.Cells.Item(60, 10).Value = 0.478
$ws.Cells.Item(60, 11).Value = -0.0000499999999999944932938
$ws.Cells.Item(60, 12).Value = $false
$ws.Cells.Item(60, 13).Value = 0.001200000000000034
$ws.Cells.Item(60, 14).Value = $false

# Row 61: 562860 生物疫苗ETF
$ws.Cells.Item(61, 1).NumberFormat = "@"
$ws.Cells.Item(61, 1).Value = "562860"
$ws.Cells.Item(61, 2).NumberFormat = "@"
$ws.Cells.Item(61, 2).Value = "2025-06-13"
$ws.Cells.Item(61, 3).Value = "生物疫苗ETF"
$ws.Cells.Item(61, 4).Value = 0.66
$ws.Cells.Item(61, 5).Value = -2.51
$ws.Cells.Item(61, 6).Value = 8.39
$ws.Cells.Item(61, 7).Value = 0.14
$ws.Cells.Item(61, 8).Value = 0.661
$ws.Cells.Item(61, 9).Value = 0.6516
$ws.Cells.Item(61, 10).Value = 0.633
$ws.Cells.Item(61, 11).Value = 0.00275000000000003
$ws.Cells.Item(61, 12).Value = $true
$ws.Cells.Item(61, 13).Value = 0.003800000000000026
$ws.Cells.Item(61, 14).Value = $false

# Row 62: 512690 酒ETF
$ws.Cells.Item(62, 1).NumberFormat = "@"
$ws.Cells.Item(62, 1).Value = "512690"
$ws.Cells.Item(62, 2).NumberFormat = "@"
$ws.Cells.Item(62, 2).Value = "2025-06-13"
$ws.Cells.Item(62, 3).Value = "酒ETF"
$ws.Cells.Item(62, 4).Value = 0.55
$ws.Cells.Item(62, 5).Value = -2.68
$ws.Cells.Item(62, 6).Value = -11.53
$ws.Cells.Item(62, 7).Value = 11.22
$ws.Cells.Item(62, 8).Value = 0.57
$ws.Cells.Item(62, 9).Value = 0.5757
$ws.Cells.Item(62, 10).Value = 0.585
$ws.Cells.Item(62, 11).Value = -0.002699999999999925
$ws.Cells.Item(62, 12).Value = $false
$ws.Cells.Item(62, 13).Value = -0.003399999999999959
$ws.Cells.Item(62, 14).Value = $false

# Row 63: 513120 港股创新药ETF
$ws.Cells.Item(63, 1).NumberFormat = "@"
$ws.Cells.Item(63, 1).Value = "513120"
$ws.Cells.Item(63, 2).NumberFormat = "@"
$ws.Cells.Item(63, 2).Value = "2025-06-13"
$ws.Cells.Item(63, 3).Value = "港股创新药ETF"
$ws.Cells.Item(63, 4).Value = 1.19
$ws.Cells.Item(63, 5).Value = -2.94
$ws.Cells.Item(63, 6).Value = 64.14
$ws.Cells.Item(63, 7).Value = 121.53
$ws.Cells.Item(63, 8).Value = 1.125
$ws.Cells.Item(63, 9).Value = 1.0739
$ws.Cells.Item(63, 10).Value = 1.006
$ws.Cells.Item(63, 11).Value = 0.01550000000000018
$ws.Cells.Item(63, 12).Value = $true
$ws.Cells.Item(63, 13).Value = 0.02100000000000013
$ws.Cells.Item(63, 14).Value = $false
